$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 and A2 hold the numeric value 0, formatted bold/boxed/centered
$b1 = $ws.Range("B1")
$b1.Value = 0
$b1.Borders.LineStyle = 1   # xlContinuous
$b1.Borders.Weight = 2      # xlThin
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4160     # xlTop

$a2 = $ws.Range("A2")
$a2.Value = 0

# Clone B1's formatting onto A2 without re-deriving the style piecemeal
# (keeps the stylesheet minimal: one new font, one new border, one new xf)
$b1.Copy()
$a2.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# B2 holds the label string, unstyled
$ws.Range("B2").Value = "disconnected_elements"
